$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.299.51"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "3.855.58"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'464.28"
$ws.Range("E5").Value = "  +9.28%  "
$ws.Range("D6").Value = "'148.46"
$ws.Range("E6").Value = "  +12.85%  "
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("E11").Value = "  -8.93%  "
$ws.Range("D12").Value = "'43.84"
$ws.Range("E12").Value = "  +6.82%  "
$ws.Range("D13").Value = "'10.42"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "4.478.41"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "'14.71"
$ws.Range("E15").Value = "  -7.43%  "
$ws.Range("D16").Value = "3.891.83"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "'20.04"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  +7.21%  "
$ws.Range("D20").Value = "67.403.00"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "'430.26"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("D22").Value = "'14.83"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "'3.27"
$ws.Range("E23").Value = "  +7.45%  "
$ws.Range("D24").Value = "'88.38"
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("E25").Value = "  +9.05%  "
$ws.Range("D26").Value = "'10.47"
$ws.Range("E26").Value = "  +14.20%  "
$ws.Range("D27").Value = "'37.54"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("D29").Value = "'5.49"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").Value = "'743.80"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'0.135"
$ws.Range("E31").Value = "  +9.79%  "
$ws.Range("D32").Value = "'13.76"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").Value = "'43.31"
$ws.Range("E34").Value = "  +11.32%  "
$ws.Range("E35").Value = "  +6.17%  "
$ws.Range("D36").Value = "'57.30"
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "'5.56"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "'0.0478"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").Value = "'0.353"
$ws.Range("E40").Value = "  +11.92%  "
$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").Value = "'2.64"
$ws.Range("E42").Value = "  +14.55%  "
$ws.Range("D43").Value = "'0.141"
$ws.Range("E43").Value = "  +4.76%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0674"
$ws.Range("E44").Value = "  -10.46%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").Value = "'3.26"
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("D48").Value = "'2.75"
$ws.Range("E48").Value = "  +7.35%  "
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'144.47"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'2.89"
$ws.Range("E51").Value = "  +2.45%  "
